$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A width matches the author's manual resize to ~24.29 "characters" (stored OOXML width).
# Excel's ColumnWidth setter quantizes to whole pixels (Maximum Digit Width = 6px for Calibri 11),
# so 23.5 is the nearest settable value that serializes to the target stored width.
$ws.Columns.Item(1).ColumnWidth = 23.5
